$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.698.14'
$ws.Range('D3').Value = '2.709.41'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.09'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '163.13'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.56%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '2.708.66'
$ws.Range('E9').Value = '  +2.32%  '
$ws.Range('E10').Value = '  -1.04%  '
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.32'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('E13').Value = '  +2.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.46'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').Value = '3.206.84'
$ws.Range('E15').Value = '  +2.42%  '
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').Value = '68.681.68'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').Value = '2.697.66'
$ws.Range('E18').Value = '  +2.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.89'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +4.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.68'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.32'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.54'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.91'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.08%  '
$ws.Range('E24').Value = '  +2.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.93'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.77%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.91'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.73%  '
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '595.78'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +6.59%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.27'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.49%  '
$ws.Range('E33').Value = '  +2.94%  '
$ws.Range('E34').Value = '  +5.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.63'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.16%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.92'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '160.74'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.76%  '
$ws.Range('E40').Value = '  +2.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.91'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.99%  '
$ws.Range('E42').Value = '  +1.57%  '
$ws.Range('E43').Value = '  +3.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.01'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.16%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('E46').Value = '  -5.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '157.80'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('E48').Value = '  +5.29%  '
$ws.Range('E49').Value = '  +5.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.607'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +6.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.04'
$ws.Range('D51').ClearFormats()
